$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.323166608810425
$ws.Range("B1").Value = 2.448666572570801
$ws.Range("C1").Value = 5.896531105041504
$ws.Range("D1").Value = 1.776394009590149
$ws.Range("E1").Value = 1.284339904785156
